$d = $word.ActiveDocument

# --- 0. Remove the pre-existing "_GoBack" bookmark from the end of the document; an
#        equivalent bookmark is (re-)created further up, inside the new "Testing Conditions"
#        paragraph, by the InsertXML call below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

function Insert-XmlBefore($range, $bodyInnerXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# --- 1. Insert "User" / "Testing Conditions" paragraphs before the "General Functionality" heading ---
$findRange = $d.Content
$findRange.Find.Execute("General Functionality") | Out-Null
$insertPoint = $d.Range($findRange.Start, $findRange.Start)

# NOTE: InsertXML merges the LAST paragraph of the inserted fragment into the target paragraph
# (adopting the target's own paragraph properties). To keep the target ("General Functionality")
# completely untouched we terminate the fragment with an empty placeholder paragraph, which is
# inserted as its own standalone (but truly empty) paragraph rather than being merged. We then
# delete that now-redundant blank paragraph afterwards.
$newParasXml = `
    '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>User</w:t></w:r></w:p>' + `
    '<w:p/>' + `
    '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Testing Conditions</w:t></w:r></w:p>' + `
    '<w:p>' + `
        '<w:r><w:t>All tests were performed on a 64-bit version of Windows 10 Professional using an 8</w:t></w:r>' + `
        '<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> generation Intel i5 CPU with 16GB of physical memory (RAM) available. The tests were all run on an x86 version of the software &#8211; except for if a System.OutOfMemoryException occurred where the test was then run again on an x64 version of the software</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> to see if the exception</w:t></w:r>' + `
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
        '<w:bookmarkEnd w:id="0"/>' + `
        '<w:r><w:t xml:space="preserve"> repeated</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve">. </w:t></w:r>' + `
    '</w:p>' + `
    '<w:p/>'

Insert-XmlBefore $insertPoint $newParasXml

# Remove the now-redundant truly-empty placeholder paragraph that sits directly before
# "General Functionality" (it exists only so the merge above did not disturb that heading).
$findRange2 = $d.Content
$findRange2.Find.Execute("General Functionality") | Out-Null
$headingStart = $findRange2.Start
$headingPara = $d.Paragraphs.Item(1)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $headingStart) {
        $headingPara = $p
        break
    }
}
$precedingPara = $headingPara.Previous()
$precedingPara.Range.Delete() | Out-Null
